$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "g"
$ws.Range("A2").Value = "hjhjkhhkh"

$ws.Range("A2").Select()
